$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Projects")
$ws.Range("C149:C158").Value = "Done"
